# Apply the vocab workbook update described by the commit:
#  - extend the "abschlieβen (abgleschlossen)" entry with " / beenden"
#  - append 12 new lesson-11 vocabulary rows (rows 294-305) to the "vocab" sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("vocab")

# ---------------------------------------------------------------------------
# 1) Extend the existing "abschlieβen (abgleschlossen)" entry (row 263, col A)
#    so that it reads "abschlieβen (abgleschlossen) / beenden".
# ---------------------------------------------------------------------------
$abschliessenCell = $ws.Range("A263")
$existingLen = $abschliessenCell.Characters().Text().Length
$appendedChars = $abschliessenCell.Characters($existingLen + 1, 0)
$appendedChars.Text = " / beenden"
$appendedRange = $abschliessenCell.Characters($existingLen + 1, 10)
$appendedRange.Font.Name = "Calibri"
$appendedRange.Font.Size = 11

# ---------------------------------------------------------------------------
# 2) Add the 12 new vocabulary rows taught on lesson date 2022-04-13 (lesson 11)
#    Column layout: A=German, B=English, C=lesson date, D=lesson number, E=type
# ---------------------------------------------------------------------------
$lessonDate = Get-Date -Year 2022 -Month 4 -Day 13 -Hour 0 -Minute 0 -Second 0
$lessonNumber = 11

# Prime the clipboard once with an existing date-formatted cell so new date
# cells reuse the same "short date" cell style instead of creating new ones.
$ws.Range("C2").Copy()
for ($r = 294; $r -le 305; $r++) {
    $ws.Range("C$r").PasteSpecial(-4122)
}

# -- German column (A) is filled in first for rows 294-301, then row 305 --
$ws.Range("A294").Value = 'das Miethaus, "-er'
$ws.Range("A295").Value = "das Mehrfamilienhaus"
$ws.Range("A296").Value = "die Eigentumswohnung"
$ws.Range("A297").Value = "die Hausverwaltung"
$ws.Range("A298").Value = "die Wohnungsagentur"
$ws.Range("A299").Value = "der Eigentumer"
$ws.Range("A300").Value = "der Burgermeistering von Paris"
$ws.Range("A301").Value = "die Innenstadt -> in der Innenstadt"
$ws.Range("A305").Value = "schicken"

# -- English column (B) is filled in next for rows 294-301 --
$ws.Range("B294").Value = "the rental house"
$ws.Range("B295").Value = "multifamily residential"
$ws.Range("B296").Value = "owner-occupied flat"
$ws.Range("B297").Value = "the property management"
$ws.Range("B298").Value = "the rental agency"
$ws.Range("B299").Value = "the owner"
$ws.Range("B300").Value = "the mayor of Paris"
$ws.Range("B301").Value = "the city center"

# -- Remaining rows 302-305 --
$ws.Range("A302").Value = "abschlieβen (abgeschlossen) / beenden"
$ws.Range("B303").Value = "to occur"
$ws.Range("B302").Value = "to finish"
$ws.Range("A303").Value = "vorliegen"
$ws.Range("A304").Value = "verschicken"
$ws.Range("B304").Value = "to send (something)"
$ws.Range("B305").Value = "to send (to someone)"

# -- Date / lesson number / type columns for all 12 new rows --
for ($r = 294; $r -le 305; $r++) {
    $ws.Range("C$r").Value = $lessonDate
    $ws.Range("D$r").Value = $lessonNumber
    $ws.Range("E$r").Value = "word"
}

# ---------------------------------------------------------------------------
# 3) Update the active selection to reflect where editing left off.
# ---------------------------------------------------------------------------
$ws.Range("A294").Select()
